$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HEATING")
$ws.Range("E35").Formula = "=1.03 * 950/1050"
